# Re-commit of changes that had accidentally only been applied inside a
# subfolder "test" copy of this workbook last time - apply them here for real.
#
# Summary of the edit:
#  - settings sheet becomes the active / selected sheet (was missingCode)
#  - "ir" sheet: two existing SKIP checkboxes get turned on (TRUE) and a new
#    trading-partner row (TPCCOV) is appended, growing the Table7 listobject
#  - "simpleton" sheet: three existing SKIP checkboxes get turned on (TRUE)
#    and a new trading-partner row (TPCCOV / COCHDQ / EXSP6 / ICS / BROK) is
#    appended, growing the Table10 listobject
#  - missingCode sheet is simply no longer the active tab

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# settings sheet — becomes the active sheet, selection moves to E3
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("settings")
[void]$wsSettings.Activate()
[void]$wsSettings.Range("E3").Select()

# ---------------------------------------------------------------------
# ir sheet — flip two SKIP checkboxes on, append TPCCOV row, grow table
# ---------------------------------------------------------------------
$wsIr = $wb.Worksheets.Item("ir")

[void]($wsIr.Range("E2").Value2 = $true)
[void]($wsIr.Range("E4").Value2 = $true)

# New row 6 inherits the formatting of row 5 (the last existing data row)
$wsIr.Range("A5:E5").Copy()
$wsIr.Range("A6:E6").PasteSpecial($xlPasteFormats)

$wsIr.Range("A6").Value2 = "TPCCOV"
$wsIr.Range("B6").Value2 = "jisavr3"
$wsIr.Range("C6").Value2 = 3
$wsIr.Range("D6").Value2 = 2

$loIr = $wsIr.ListObjects.Item(1)
[void]$loIr.Resize($wsIr.Range("A1:E6"))

[void]$wsIr.Range("A7").Select()

# ---------------------------------------------------------------------
# simpleton sheet — flip three SKIP checkboxes on, append TPCCOV row,
# grow table
# ---------------------------------------------------------------------
$wsSimpleton = $wb.Worksheets.Item("simpleton")

[void]($wsSimpleton.Range("H2").Value2 = $true)
[void]($wsSimpleton.Range("H3").Value2 = $true)
[void]($wsSimpleton.Range("H5").Value2 = $true)

# New row 6 formatting: first two columns copy from row 5 directly, the
# middle columns (C:G) pick up the plain "below table" border style that
# already lives on G5, and H copies row 5's right-edge style.
$wsSimpleton.Range("A5:B5").Copy()
$wsSimpleton.Range("A6:B6").PasteSpecial($xlPasteFormats)

$wsSimpleton.Range("G5").Copy()
$wsSimpleton.Range("C6:G6").PasteSpecial($xlPasteFormats)

$wsSimpleton.Range("H5").Copy()
$wsSimpleton.Range("H6").PasteSpecial($xlPasteFormats)

$wsSimpleton.Range("A6").Value2 = "TPCCOV"
$wsSimpleton.Range("C6").Value2 = "COCHDQ"
$wsSimpleton.Range("D6").Value2 = "EXSP6"
$wsSimpleton.Range("E6").Value2 = "ICS"
$wsSimpleton.Range("F6").Value2 = "BROK"

$loSimpleton = $wsSimpleton.ListObjects.Item(1)
[void]$loSimpleton.Resize($wsSimpleton.Range("A1:H6"))

[void]$wsSimpleton.Range("G6").Select()

# ---------------------------------------------------------------------
# missingCode sheet — no longer the active tab; selection unaffected
# ---------------------------------------------------------------------
$wsMissingCode = $wb.Worksheets.Item("missingCode")
[void]$wsMissingCode.Range("D12").Select()

# settings is the sheet that should end up active/selected
[void]$wsSettings.Activate()
